$d = $word.ActiveDocument

# --- Change 1: rename the "chapter-one" bookmark/heading to
#     "discrete-versus-continuous" / "Discrete versus continuous".
#
# This COM host does not wire up Bookmark.Delete (and Bookmark.Name is
# read-only here), so a plain rename-in-place isn't available. Instead we
# insert a fresh Heading1 paragraph carrying the new text + bookmark right
# before the old one, then delete the old paragraph outright -- this keeps
# the bookmark table free of duplicates/leftovers.
$oldHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq "Chapter One") {
        $oldHeading = $para
        break
    }
}

$oldIndex = $oldHeading.Index
$oldHeading.Range.InsertParagraphBefore()
$newHeading = $d.Paragraphs.Item($oldIndex)
$newHeading.Range.Text = "Discrete versus continuous"

$bmRange = $d.Paragraphs.Item($oldIndex).Range
$bmRange.MoveEnd(1, -1)
$d.Bookmarks.Add("discrete-versus-continuous", $bmRange)

$d.Paragraphs.Item($oldIndex + 1).Range.Delete()

# --- Change 2: add a new "Something from Vimeo:" BodyText paragraph right
#     after the "Strikeout" paragraph.
$strikePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq "Strikeout") {
        $strikePara = $para
        break
    }
}

$strikeIndex = $strikePara.Index
$strikePara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($strikeIndex + 1)
$newPara.Range.Text = "Something from Vimeo:"
$newPara.Style = "BodyText"
